$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Formula = "'001"
$ws.Range("M2").Value = "2020-12-17 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 18233340.34
$ws.Range("P2").Value = 1861611886.82
$ws.Range("Q2").Value = 1840335733.27
$ws.Range("S2").Value = 1680814621.66
$ws.Range("T2").Value = 1680814621.66
$ws.Range("V2").Value = 19448585.36
$ws.Range("W2").Value = 32049881.2
$ws.Range("X2").Value = 35327673.89
$ws.Range("Y2").Value = 25135946.81
$ws.Range("Z2").Value = 25492619.02
$ws.Range("AA2").Value = 7259278.68
$ws.Range("AG2").Value = 3372373.04
$ws.Range("AS2").Value = 11700680.43
